# Update data/example to newest format
# Fix missing trailing / if using ExcelWriter from command line
#
# Set the title row (row 1) on every worksheet to a fixed, custom height
# of 24 points (matches the newer ExcelWriter output format).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Rows.Item(1).RowHeight = 24
}
